# Update attendance/sales figures in column F for both the "展览" and
# "全部类型" worksheets (which hold duplicated data).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# "展览" sheet updates
$wsExhibit.Range("F5").Value  = 3073
$wsExhibit.Range("F7").Value  = 2334
$wsExhibit.Range("F9").Value  = 115
$wsExhibit.Range("F11").Value = 1123
$wsExhibit.Range("F15").Value = 722
$wsExhibit.Range("F16").Value = 278
$wsExhibit.Range("F17").Value = 300
$wsExhibit.Range("F18").Value = 13
$wsExhibit.Range("F23").Value = 7

# "全部类型" sheet updates
$wsAll.Range("F5").Value  = 3073
$wsAll.Range("F7").Value  = 2334
$wsAll.Range("F9").Value  = 115
$wsAll.Range("F11").Value = 1123
$wsAll.Range("F15").Value = 723
$wsAll.Range("F16").Value = 278
$wsAll.Range("F17").Value = 300
$wsAll.Range("F18").Value = 13
$wsAll.Range("F23").Value = 7
